$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content fix: the "ACT_DNS" rows under the Ramp-Down block were mis-tagged
#     with the Ramp-Up attribute's text ("ACT_DNS" survives nowhere else in the
#     workbook); correct B6:B9 to read "ACT_UPS" like the rest of the table.
$ws.Range("B6").Value = "ACT_UPS"
$ws.Range("B7").Value = "ACT_UPS"
$ws.Range("B8").Value = "ACT_UPS"
$ws.Range("B9").Value = "ACT_UPS"

# B6 had no explicit formatting (unlike its row neighbours C6/D6); bring it into
# line with the rest of row 6 by copying that formatting across.
$ws.Range("C6").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Reposition / resize the second picture (moved further right/down and
#     slightly resized) ---
$pic = $ws.Shapes.Item(2)
$pic.Left = 637.884094488189
$pic.Top = 252.59992125984252
$pic.Width = 447.93149606299215
$pic.Height = 316.4993700787402

# --- View state: zoom in to 90% and move the selection to F8 ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 90
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F8").Select()
